$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Hunk 0: ALC row 19
$ws_ALC.Range("H19").Value = 642.8570999999999
$ws_ALC.Range("I19").Value = 800.5
$ws_ALC.Range("J19").Value = 579.8
$ws_ALC.Range("K19").Value = 800.5
$ws_ALC.Range("L19").Value = 579.8
$ws_ALC.Range("M19").Value = -625.5
$ws_ALC.Range("N19").Value = -929.8

# Hunk 1: ALC row 53
$ws_ALC.Range("H53").Value = 396.08334
$ws_ALC.Range("I53").Value = 73.666664
$ws_ALC.Range("K53").Value = 73.666664
$ws_ALC.Range("M53").Value = 563.333336

# Hunk 2: ALC row 125
$ws_ALC.Range("H125").Value = 1910
$ws_ALC.Range("J125").Value = 2046.6666
$ws_ALC.Range("L125").Value = 18419.9994
$ws_ALC.Range("N125").Value = -23339.9994

# Hunk 3: ALC row 130
$ws_ALC.Range("H130").Value = 40780
$ws_ALC.Range("I130").Value = 0
$ws_ALC.Range("J130").Value = 40780
$ws_ALC.Range("K130").Value = 0
$ws_ALC.Range("L130").Value = 40780
$ws_ALC.Range("N130").Value = -50820
$ws_ALC.Range("M130").ClearContents()

# Hunk 4: ALC row 132
$ws_ALC.Range("H132").Value = 3476461.5
$ws_ALC.Range("I132").Value = 3575503.2
$ws_ALC.Range("J132").Value = 10000
$ws_ALC.Range("K132").Value = 10726509.6
$ws_ALC.Range("L132").Value = 30000
$ws_ALC.Range("M132").Value = -10723979.6
$ws_ALC.Range("N132").Value = -35060

# Hunk 5: ALC row 133
$ws_ALC.Range("H133").Value = 66662.22
$ws_ALC.Range("J133").Value = 66662.22
$ws_ALC.Range("L133").Value = 66662.22
$ws_ALC.Range("N133").Value = -76782.22

# Hunk 6: ARM row 32
$ws_ARM.Range("H32").Value = 5182.0576
$ws_ARM.Range("I32").Value = 3651.2046
$ws_ARM.Range("J32").Value = 13601.75
$ws_ARM.Range("K32").Value = 3651.2046
$ws_ARM.Range("L32").Value = 13601.75
$ws_ARM.Range("M32").Value = -3364.2046
$ws_ARM.Range("N32").Value = -14175.75

# Hunk 7: ARM row 61
$ws_ARM.Range("H61").Value = 1524.5555
$ws_ARM.Range("I61").Value = 1361.3846
$ws_ARM.Range("J61").Value = 1948.8
$ws_ARM.Range("K61").Value = 1361.3846
$ws_ARM.Range("L61").Value = 1948.8
$ws_ARM.Range("M61").Value = -1149.3846
$ws_ARM.Range("N61").Value = -2372.8

# Hunk 8: ARM row 132
$ws_ARM.Range("H132").Value = 27624.75
$ws_ARM.Range("I132").Value = 35166.668
$ws_ARM.Range("J132").Value = 4999
$ws_ARM.Range("K132").Value = 105500.004
$ws_ARM.Range("L132").Value = 14997
$ws_ARM.Range("M132").Value = -102970.004
$ws_ARM.Range("N132").Value = -20057

# Hunk 9: ARM row 136
$ws_ARM.Range("H136").Value = 1524.5555
$ws_ARM.Range("I136").Value = 1361.3846
$ws_ARM.Range("J136").Value = 1948.8
$ws_ARM.Range("K136").Value = 4084.1538
$ws_ARM.Range("L136").Value = 5846.4
$ws_ARM.Range("M136").Value = -1534.1538
$ws_ARM.Range("N136").Value = -10946.4

# Hunk 10: BSM row 22
$ws_BSM.Range("H22").Value = 337.36365
$ws_BSM.Range("I22").Value = 353.5
$ws_BSM.Range("J22").Value = 318
$ws_BSM.Range("K22").Value = 353.5
$ws_BSM.Range("L22").Value = 318
$ws_BSM.Range("M22").Value = -180.5
$ws_BSM.Range("N22").Value = -664

# Hunk 11: BSM row 81
$ws_BSM.Range("H81").Value = 26353.334
$ws_BSM.Range("J81").Value = 26353.334
$ws_BSM.Range("L81").Value = 26353.334
$ws_BSM.Range("N81").Value = -28475.334

# Hunk 12: BSM row 84
$ws_BSM.Range("H84").Value = 26353.334
$ws_BSM.Range("J84").Value = 26353.334
$ws_BSM.Range("L84").Value = 79060.00199999999
$ws_BSM.Range("N84").Value = -89668.00199999999

# Hunk 13: BSM row 94
$ws_BSM.Range("H94").Value = 1129.8823
$ws_BSM.Range("I94").Value = 1120.5333
$ws_BSM.Range("J94").Value = 1200
$ws_BSM.Range("K94").Value = 1120.5333
$ws_BSM.Range("L94").Value = 1200
$ws_BSM.Range("M94").Value = -669.5333000000001
$ws_BSM.Range("N94").Value = -2102

# Hunk 14: BSM row 132
$ws_BSM.Range("H132").Value = 45780
$ws_BSM.Range("J132").Value = 45780
$ws_BSM.Range("L132").Value = 45780
$ws_BSM.Range("N132").Value = -55900

# Hunk 15: BSM row 138
$ws_BSM.Range("H138").Value = 81152
$ws_BSM.Range("J138").Value = 81152
$ws_BSM.Range("L138").Value = 81152
$ws_BSM.Range("N138").Value = -91432

# Hunk 16: CRP row 22
$ws_CRP.Range("H22").Value = 1041.5385
$ws_CRP.Range("I22").Value = 1351.1111
$ws_CRP.Range("J22").Value = 345
$ws_CRP.Range("K22").Value = 1351.1111
$ws_CRP.Range("L22").Value = 345
$ws_CRP.Range("M22").Value = -1001.1111
$ws_CRP.Range("N22").Value = -1045

# Hunk 17: CRP row 130
$ws_CRP.Range("H130").Value = 39566.668
$ws_CRP.Range("J130").Value = 39566.668
$ws_CRP.Range("L130").Value = 39566.668
$ws_CRP.Range("N130").Value = -49606.668

# Hunk 18: CRP row 132
$ws_CRP.Range("H132").Value = 4341.516
$ws_CRP.Range("I132").Value = 799.5925999999999
$ws_CRP.Range("J132").Value = 28249.5
$ws_CRP.Range("K132").Value = 2398.7778
$ws_CRP.Range("L132").Value = 84748.5
$ws_CRP.Range("M132").Value = 131.2222000000002
$ws_CRP.Range("N132").Value = -89808.5

# Hunk 19: CUL row 12
$ws_CUL.Range("H12").Value = 57.615383
$ws_CUL.Range("I12").Value = 92.333336
$ws_CUL.Range("J12").Value = 47.2
$ws_CUL.Range("K12").Value = 277.000008
$ws_CUL.Range("L12").Value = 141.6
$ws_CUL.Range("M12").Value = -104.000008
$ws_CUL.Range("N12").Value = -487.6

# Hunk 20: CUL row 25
$ws_CUL.Range("H25").Value = 2650
$ws_CUL.Range("I25").Value = 1800
$ws_CUL.Range("J25").Value = 3500
$ws_CUL.Range("K25").Value = 5400
$ws_CUL.Range("L25").Value = 10500
$ws_CUL.Range("M25").Value = -5231
$ws_CUL.Range("N25").Value = -10838

# Hunk 21: CUL row 30
$ws_CUL.Range("H30").Value = 2650
$ws_CUL.Range("I30").Value = 1800
$ws_CUL.Range("J30").Value = 3500
$ws_CUL.Range("K30").Value = 5400
$ws_CUL.Range("L30").Value = 10500
$ws_CUL.Range("M30").Value = -5298
$ws_CUL.Range("N30").Value = -10704

# Hunk 22: GSM row 2
$ws_GSM.Range("H2").Value = 128.33333
$ws_GSM.Range("I2").Value = 57.5
$ws_GSM.Range("K2").Value = 57.5
$ws_GSM.Range("M2").Value = 55.5

# Hunk 23: GSM row 25
$ws_GSM.Range("H25").Value = 1800
$ws_GSM.Range("J25").Value = 1800
$ws_GSM.Range("L25").Value = 1800
$ws_GSM.Range("N25").Value = -2858

# Hunk 24: GSM row 122
$ws_GSM.Range("H122").Value = 2882.5
$ws_GSM.Range("I122").Value = 2853.5
$ws_GSM.Range("J122").Value = 2897
$ws_GSM.Range("K122").Value = 8560.5
$ws_GSM.Range("L122").Value = 8691
$ws_GSM.Range("M122").Value = -6110.5
$ws_GSM.Range("N122").Value = -13591

# Hunk 25: GSM row 126
$ws_GSM.Range("H126").Value = 144033.14
$ws_GSM.Range("I126").Value = 250948
$ws_GSM.Range("J126").Value = 1480
$ws_GSM.Range("K126").Value = 752844
$ws_GSM.Range("L126").Value = 4440
$ws_GSM.Range("M126").Value = -750374
$ws_GSM.Range("N126").Value = -9380

# Hunk 26: GSM row 133
$ws_GSM.Range("H133").Value = 50780
$ws_GSM.Range("J133").Value = 50780
$ws_GSM.Range("L133").Value = 50780
$ws_GSM.Range("N133").Value = -60900

# Hunk 27: LTW row 16
$ws_LTW.Range("H16").Value = 864
$ws_LTW.Range("I16").Value = 300.5
$ws_LTW.Range("J16").Value = 2366.6667
$ws_LTW.Range("K16").Value = 300.5
$ws_LTW.Range("L16").Value = 2366.6667
$ws_LTW.Range("M16").Value = -130.5
$ws_LTW.Range("N16").Value = -2706.6667

# Hunk 28: LTW row 22
$ws_LTW.Range("H22").Value = 683.0952
$ws_LTW.Range("I22").Value = 564.0909
$ws_LTW.Range("J22").Value = 814
$ws_LTW.Range("K22").Value = 564.0909
$ws_LTW.Range("L22").Value = 814
$ws_LTW.Range("M22").Value = -269.0909
$ws_LTW.Range("N22").Value = -1404

# Hunk 29: LTW row 27
$ws_LTW.Range("H27").Value = 683.0952
$ws_LTW.Range("I27").Value = 564.0909
$ws_LTW.Range("J27").Value = 814
$ws_LTW.Range("K27").Value = 564.0909
$ws_LTW.Range("L27").Value = 814
$ws_LTW.Range("M27").Value = -457.0909
$ws_LTW.Range("N27").Value = -1028

# Hunk 30: LTW row 46
$ws_LTW.Range("H46").Value = 800
$ws_LTW.Range("I46").Value = 775
$ws_LTW.Range("J46").Value = 1000
$ws_LTW.Range("K46").Value = 775
$ws_LTW.Range("L46").Value = 1000
$ws_LTW.Range("M46").Value = -587
$ws_LTW.Range("N46").Value = -1376

# Hunk 31: LTW row 60
$ws_LTW.Range("H60").Value = 0
$ws_LTW.Range("J60").Value = 0
$ws_LTW.Range("L60").Value = 0
$ws_LTW.Range("N60").ClearContents()

# Hunk 32: LTW row 82
$ws_LTW.Range("H82").Value = 2466.238
$ws_LTW.Range("I82").Value = 1360
$ws_LTW.Range("J82").Value = 2811.9375
$ws_LTW.Range("K82").Value = 1360
$ws_LTW.Range("L82").Value = 2811.9375
$ws_LTW.Range("M82").Value = -999
$ws_LTW.Range("N82").Value = -3533.9375

# Hunk 33: LTW row 85
$ws_LTW.Range("H85").Value = 2466.238
$ws_LTW.Range("I85").Value = 1360
$ws_LTW.Range("J85").Value = 2811.9375
$ws_LTW.Range("K85").Value = 1360
$ws_LTW.Range("L85").Value = 2811.9375
$ws_LTW.Range("M85").Value = -112
$ws_LTW.Range("N85").Value = -5307.9375

# Hunk 34: LTW row 100
$ws_LTW.Range("H100").Value = 2760.6128
$ws_LTW.Range("I100").Value = 1942.7142
$ws_LTW.Range("K100").Value = 1942.7142
$ws_LTW.Range("M100").Value = -1401.7142

# Hunk 35: LTW row 114
$ws_LTW.Range("H114").Value = 40398
$ws_LTW.Range("J114").Value = 40398
$ws_LTW.Range("L114").Value = 40398
$ws_LTW.Range("N114").Value = -49076

# Hunk 36: LTW row 122
$ws_LTW.Range("H122").Value = 4907.6665
$ws_LTW.Range("I122").Value = 6122.222
$ws_LTW.Range("J122").Value = 3085.8333
$ws_LTW.Range("K122").Value = 18366.666
$ws_LTW.Range("L122").Value = 9257.499899999999
$ws_LTW.Range("M122").Value = -15916.666
$ws_LTW.Range("N122").Value = -14157.4999

# Hunk 37: LTW row 125
$ws_LTW.Range("H125").Value = 40643
$ws_LTW.Range("J125").Value = 40643
$ws_LTW.Range("L125").Value = 40643
$ws_LTW.Range("N125").Value = -50483

# Hunk 38: LTW row 130
$ws_LTW.Range("H130").Value = 35194.5
$ws_LTW.Range("J130").Value = 35194.5
$ws_LTW.Range("L130").Value = 35194.5
$ws_LTW.Range("N130").Value = -45234.5

# Hunk 39: LTW row 132
$ws_LTW.Range("H132").Value = 4964.1064
$ws_LTW.Range("I132").Value = 4810.814
$ws_LTW.Range("J132").Value = 6612
$ws_LTW.Range("K132").Value = 14432.442
$ws_LTW.Range("L132").Value = 19836
$ws_LTW.Range("M132").Value = -11902.442
$ws_LTW.Range("N132").Value = -24896

# Hunk 40: WVR row 21
$ws_WVR.Range("H21").Value = 2862999.8
$ws_WVR.Range("I21").Value = 10005000
$ws_WVR.Range("J21").Value = 6200
$ws_WVR.Range("K21").Value = 10005000
$ws_WVR.Range("L21").Value = 6200
$ws_WVR.Range("M21").Value = -10004765
$ws_WVR.Range("N21").Value = -6670

# Hunk 41: WVR row 35
$ws_WVR.Range("H35").Value = 2862999.8
$ws_WVR.Range("I35").Value = 10005000
$ws_WVR.Range("J35").Value = 6200
$ws_WVR.Range("K35").Value = 10005000
$ws_WVR.Range("L35").Value = 6200
$ws_WVR.Range("M35").Value = -10004710
$ws_WVR.Range("N35").Value = -6780

# Hunk 42: WVR row 107
$ws_WVR.Range("H107").Value = 400
$ws_WVR.Range("I107").Value = 400
$ws_WVR.Range("J107").Value = 0
$ws_WVR.Range("K107").Value = 1200
$ws_WVR.Range("L107").Value = 0
$ws_WVR.Range("M107").Value = 720
$ws_WVR.Range("N107").ClearContents()

Write-Output "Applied all edits"
